$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @(2, 0.6666666666666666, 0.262426, 0.7872779999999999, 0.2164351829578579, 0.2164351829578578, 1, 0.3333333333333333, 0.2972546666666667, 0.891764, 0.02381895177797704, 0.02381895177797704, 0.07800735315466667, 0.7020661783919999, 0.005155259185930855, 0.005155259185930854)
$data[3] = @(2, 0.6666666666666666, 0.262426, 0.7872779999999999, 0.2164351829578579, 0.2164351829578578, 3, 1, 9.750671000000001, 29.252013, 0.7813191461594744, 0.7813191461594745, 2.558829587846, 23.029466290614, 0.1691049523475031, 0.1691049523475031)
$data[4] = @(2, 0.6666666666666666, 0.262426, 0.7872779999999999, 0.2164351829578579, 0.2164351829578578, 3, 1, 2.431828666666667, 7.295485999999999, 0.1948619020625486, 0.1948619020625486, 0.6381750696786667, 5.743575627107999, 0.04217497142442389, 0.04217497142442388)
$data[5] = @(2, 0.6666666666666666, 0.430825, 1.292475, 0.3553218343373718, 0.3553218343373717, 1, 0.3333333333333333, 0.2972546666666667, 0.891764, 0.02381895177797704, 0.02381895177797704, 0.1280647417666667, 1.1525826759, 0.008463393637744204, 0.008463393637744203)
$data[6] = @(2, 0.6666666666666666, 0.430825, 1.292475, 0.3553218343373718, 0.3553218343373717, 3, 1, 9.750671000000001, 29.252013, 0.7813191461594744, 0.7813191461594745, 4.200832833575, 37.807495502175, 0.2776197522162935, 0.2776197522162935)
$data[7] = @(2, 0.6666666666666666, 0.430825, 1.292475, 0.3553218343373718, 0.3553218343373717, 3, 1, 2.431828666666667, 7.295485999999999, 0.1948619020625486, 0.1948619020625486, 1.047692585316667, 9.42923326785, 0.06923868848333406, 0.06923868848333405)
$data[8] = @(1, 0.3333333333333333, 0.1998913333333333, 0.599674, 0.1648598740280695, 0.1648598740280695, 1, 0.3333333333333333, 0.2972546666666667, 0.891764, 0.02381895177797704, 0.02381895177797704, 0.05941863165955556, 0.534767684936, 0.003926789389597955, 0.003926789389597956)
$data[9] = @(1, 0.3333333333333333, 0.1998913333333333, 0.599674, 0.1648598740280695, 0.1648598740280695, 3, 1, 9.750671000000001, 29.252013, 0.7813191461594744, 0.7813191461594745, 1.949074627084667, 17.541671643762, 0.1288081760115697, 0.1288081760115697)
$data[10] = @(1, 0.3333333333333333, 0.1998913333333333, 0.599674, 0.1648598740280695, 0.1648598740280695, 3, 1, 2.431828666666667, 7.295485999999999, 0.1948619020625486, 0.1948619020625486, 0.4861014746182222, 4.374913271564, 0.03212490862690177, 0.03212490862690177)
$data[11] = @(2, 0.6666666666666666, 0.31935, 0.95805, 0.2633831086767009, 0.2633831086767008, 1, 0.3333333333333333, 0.2972546666666667, 0.891764, 0.02381895177797704, 0.02381895177797704, 0.09492827779999999, 0.8543545001999999, 0.006273509564704024, 0.006273509564704024)
$data[12] = @(2, 0.6666666666666666, 0.31935, 0.95805, 0.2633831086767009, 0.2633831086767008, 3, 1, 9.750671000000001, 29.252013, 0.7813191461594744, 0.7813191461594745, 3.11387678385, 28.02489105465, 0.205786265584108, 0.205786265584108)
$data[13] = @(2, 0.6666666666666666, 0.31935, 0.95805, 0.2633831086767009, 0.2633831086767008, 3, 1, 2.431828666666667, 7.295485999999999, 0.1948619020625486, 0.1948619020625486, 0.7766044846999999, 6.989440362299999, 0.05132333352788888, 0.05132333352788888)

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in 2..13) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$r"
        $ws.Range($cellRef).Value = $rowVals[$i]
    }
}
